$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (D1:W1): "Day-Hour" -> "DayHour" (drop the dash) ---
$headers = @{
    "D1" = "Mon08";  "E1" = "Mon010"; "F1" = "Mon012"; "G1" = "Mon02"
    "H1" = "Tue08";  "I1" = "Tue010"; "J1" = "Tue012"; "K1" = "Tue02"
    "L1" = "Wed08";  "M1" = "Wed010"; "N1" = "Wed012"; "O1" = "Wed02"
    "P1" = "Thu08";  "Q1" = "Thu010"; "R1" = "Thu012"; "S1" = "Thu02"
    "T1" = "Fri08";  "U1" = "Fri010"; "V1" = "Fri012"; "W1" = "Fri02"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# --- Replace the placeholder "-" spot-count cells with numeric 0 ---
$zeroCells = @("D6","J6","N6","O6","R6","D17","D25","D26","D31","D32","D33","D34")
foreach ($addr in $zeroCells) {
    $ws.Range($addr).Value = 0
}

# --- Update the view: scroll back to A1 (clear the B1 top-left freeze) and move the selection ---
[void]$ws.Range("A1").Select()
[void]$ws.Range("J12").Select()
